$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.460.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.40%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.896.49"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.96%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.06%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'238.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.02%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  -0.06%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.95%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.2922"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.54%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.06684"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.31%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'1.907.39"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.54%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'16.93"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.91%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07336"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.42%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'5.175"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +3.54%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'87.54"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.46%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.6659"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +2.18%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'30.443.93"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.28%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'13.46"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +3.94%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.000007828"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.42%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'2.166.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +2.16%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'5.320"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +13.02%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  +0.02%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'192.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.44%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'6.113"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.20%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'9.486"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +2.24%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'162.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +3.10%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -0.81%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'1.936"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +6.16%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'1.471"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +4.75%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'4.326"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +2.07%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.09160"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +1.79%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'4.052"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +3.38%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.05172"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.91%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'0.7391"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +2.43%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.102"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +2.25%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'2.719"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +0.97%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -0.21%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'2.677"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.67%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.9235"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.76%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  -0.33%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.4383"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.08%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'106.86"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'5.906"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +3.49%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -0.25%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'68.77"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +20.73%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  +2.73%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'7.579"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +3.20%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'8.984"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +3.73%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'34.84"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +5.43%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.05854"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'0.3919"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -2.14%  "
$ws.Range("E51").Style = "Normal"
